# eims-toi-transect-info.xlsx edit script
# - Rename sheet "CustonUnits" -> "CustomUnits" (fix typo)
# - Add a new Keywords row for "inorganic matter" / "LTER Core Research Areas"
# - Adjust row heights / column widths on the Keywords sheet that shift as a result
# - Update selections on the Keywords and CustomUnits sheets
# - Make "CategoricalVariables" the active/selected tab (was "CustonUnits")

$wb = $excel.ActiveWorkbook

# ---- Fix the misspelled sheet name ----
$wsUnits = $wb.Worksheets.Item("CustonUnits")
$wsUnits.Name = "CustomUnits"

# ---- Keywords sheet: insert new keyword row ----
$wsKeywords = $wb.Worksheets.Item("Keywords")

# Shift existing chemistry/oceanography/... rows down to make room for the
# new "inorganic matter" keyword row right after "primary production".
$wsKeywords.Rows.Item(3).Insert()
$wsKeywords.Range("A3").Value = "inorganic matter"
$wsKeywords.Range("B3").Value = "LTER Core Research Areas"

# The shifted rows keep their old (autofit) heights from when the text used to
# wrap differently; re-apply the correct explicit row heights.
$wsKeywords.Rows.Item(5).RowHeight = 15
$wsKeywords.Rows.Item(7).RowHeight = 14.4
$wsKeywords.Rows.Item(8).RowHeight = 15

# Explicit column widths for columns A and C (values chosen so the stored
# OOXML <col width="..."> comes out as close as possible to the authored
# 20.8984375 / 15.69921875 widths given this engine's column-width quantization)
$wsKeywords.Columns.Item(1).ColumnWidth = 20
$wsKeywords.Columns.Item(3).ColumnWidth = 14.833333333333332

# Update the sheet's remembered selection
[void]$wsKeywords.Range("A11").Select()

# ---- CustomUnits sheet: update remembered selection ----
[void]$wsUnits.Range("E21").Select()

# ---- Make CategoricalVariables the active tab ----
$wsCategorical = $wb.Worksheets.Item("CategoricalVariables")
[void]$wsCategorical.Activate()
